# Trade #5 closed at 2026-02-17 13:33:39 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("Summary")
$strategyStatus = $wb.Worksheets.Item("Strategy Status")
$allTrades = $wb.Worksheets.Item("All Trades")
$marketMaking = $wb.Worksheets.Item("MarketMaking")

# --- Summary sheet updates ---
$summary.Range("B3").Value = 1199.92
$summary.Range("B4").Value = -0.08
$summary.Range("B5").Value = -0.32
$summary.Range("B6").Value = 5
$summary.Range("B8").Value = 4
$summary.Range("B9").Value = 20

# --- Strategy Status sheet updates (MarketMaking row) ---
$strategyStatus.Range("C4").Value = 99.92
$strategyStatus.Range("D4").Value = 5
$strategyStatus.Range("E4").Value = -0.08
$strategyStatus.Range("F4").Value = -0.08
$strategyStatus.Range("G4").Value = 20

# --- New trade row (Trade #5) appended to both "All Trades" and "MarketMaking" sheets ---
foreach ($sheet in @($allTrades, $marketMaking)) {
    $sheet.Range("A6").Value = 5
    # Leading apostrophe forces the date-shaped string to stay text instead
    # of being auto-converted to a date serial number.
    $sheet.Range("B6").Value = "'2026-02-17"
    $sheet.Range("C6").Value = "13:33:32"
    $sheet.Range("D6").Value = "MarketMaking"
    $sheet.Range("E6").Value = "UP"
    $sheet.Range("F6").Value = 0.06
    $sheet.Range("G6").Value = 0.05
    $sheet.Range("H6").Value = "CLOSED"
    $sheet.Range("I6").Value = -16.6667
    $sheet.Range("J6").Value = -0.01
    $sheet.Range("K6").Value = 99.92
    $sheet.Range("L6").Value = 0
    $sheet.Range("M6").Value = 0
    $sheet.Range("N6").Value = 0.6
    $sheet.Range("O6").Value = "Normal spread capture: 19600 bps"
    $sheet.Range("P6").Value = "early_exit"
    $sheet.Range("Q6").Value = 0.13
}
